# Add "Wins", "Losses", "Ties" columns (AD, AE, AF) with season record data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used data row (excluding header)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Copy the header formatting (bold, bordered, centered) from an existing
# header cell (A1) onto the three new header cells, then set their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill season record values for each data row (row 2 through last row)
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 74   # AD = Wins
    $ws.Cells.Item($r, 31).Value = 88   # AE = Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF = Ties
}
